$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Plain text / URL / label updates (unambiguous as text) ---
$ws.Range("D2").Value = "43.761.47"
$ws.Range("E2").Value = "  +4.53%  "
$ws.Range("D3").Value = "2.281.76"
$ws.Range("E3").Value = "  +2.20%  "
$ws.Range("E4").Value = "  +0.01%  "
$ws.Range("E5").Value = "  -0.51%  "
$ws.Range("E6").Value = "  +0.82%  "
$ws.Range("E7").Value = "  +6.04%  "
$ws.Range("E8").Value = "  +0.01%  "
$ws.Range("E9").Value = "  +4.85%  "
$ws.Range("E10").Value = "  +5.31%  "
$ws.Range("E11").Value = "  -0.72%  "
$ws.Range("E12").Value = "  +17.58%  "
$ws.Range("E13").Value = "  +0.17%  "
$ws.Range("D14").Value = "2.620.34"
$ws.Range("E14").Value = "  +2.30%  "
$ws.Range("E15").Value = "  +0.33%  "
$ws.Range("E16").Value = "  +4.85%  "
$ws.Range("E17").Value = "  +1.58%  "
$ws.Range("D18").Value = "2.298.06"
$ws.Range("E18").Value = "  +2.66%  "
$ws.Range("D19").Value = "43.718.40"
$ws.Range("E19").Value = "  +4.73%  "
$ws.Range("D20").Value = "0.0₃0947"
$ws.Range("E20").Value = "  +4.11%  "
$ws.Range("E21").Value = "  +1.12%  "
$ws.Range("E22").Value = "  +0.78%  "
$ws.Range("E23").Value = "  +1.23%  "
$ws.Range("E24").Value = "  +0.11%  "
$ws.Range("E25").Value = "  +7.50%  "
$ws.Range("E26").Value = "  -1.51%  "
$ws.Range("E27").Value = "  +2.01%  "
$ws.Range("E28").Value = "  +1.41%  "
$ws.Range("E29").Value = "  -2.88%  "
$ws.Range("E30").Value = "  +2.87%  "
$ws.Range("E31").Value = "  +2.86%  "
$ws.Range("E32").Value = "  +4.17%  "
$ws.Range("E33").Value = "  +0.25%  "
$ws.Range("E34").Value = "  +6.18%  "
$ws.Range("E35").Value = "  +2.18%  "
$ws.Range("E36").Value = "  +0.46%  "
$ws.Range("E37").Value = "  +0.97%  "
$ws.Range("E38").Value = "  +2.52%  "
$ws.Range("E39").Value = "  -2.00%  "
$ws.Range("E40").Value = "  +3.21%  "
$ws.Range("E41").Value = "  -0.07%  "
$ws.Range("E42").Value = "  +27.97%  "
$ws.Range("E43").Value = "  +3.33%  "
$ws.Range("B44").Value = "TerraClassic"
$ws.Range("C44").Value = "https://coinranking.com/coin/AaQUAs2Mc+terraclassic-lunc"
$ws.Range("E44").Value = "  -5.48%  "
$ws.Range("B45").Value = "FraxShare"
$ws.Range("C45").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("E45").Value = "  -1.67%  "
$ws.Range("E46").Value = "  -0.64%  "
$ws.Range("E47").Value = "  +0.28%  "
$ws.Range("E48").Value = "  -0.76%  "
$ws.Range("D49").Value = "1.484.03"
$ws.Range("E49").Value = "  +0.86%  "
$ws.Range("E50").Value = "  +1.49%  "
$ws.Range("E51").Value = "  +3.06%  "

# --- Numeric-looking price values that must remain stored as text ---
# (the workbook stores prices as plain text strings, e.g. "8.47", not numbers,
#  so we force text format, assign, then restore the default "Normal" style
#  to avoid leaving a custom number format applied to the cell)
$textPrices = [ordered]@{
    "D5" = "231.53"
    "D7" = "64.48"
    "D9" = "0.426"
    "D10" = "0.0954"
    "D11" = "57.66"
    "D12" = "26.71"
    "D15" = "15.74"
    "D17" = "0.815"
    "D21" = "73.36"
    "D22" = "6.16"
    "D23" = "250.74"
    "D25" = "2.57"
    "D27" = "9.88"
    "D28" = "171.60"
    "D29" = "0.138"
    "D30" = "20.52"
    "D32" = "2.75"
    "D34" = "0.0699"
    "D37" = "6.63"
    "D39" = "2.35"
    "D42" = "11.06"
    "D43" = "4.64"
    "D44" = "0.000223"
    "D45" = "8.47"
    "D46" = "1.22"
    "D47" = "0.0965"
    "D48" = "97.99"
    "D50" = "16.86"
    "D51" = "2.34"
}

foreach ($addr in $textPrices.Keys) {
    $cell = $ws.Range($addr)
    $cell.NumberFormat = "@"
    $cell.Value = $textPrices[$addr]
    $cell.Style = "Normal"
}
